# Add a new enquiry row (row 8) to the enquiry list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A8: ID, B8: Camp ID
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 2

# C8: Details - typed, then edited (matches how the entry was authored: an
# initial draft of the enquiry text was typed in and then corrected).
$ws.Range("C8").Value = "This is a sample enquiry"
$ws.Range("C8").Value = "This is an edited sample enquiry"

# D8: Answer - left blank (no answer yet for this enquiry).
$ws.Range("D8").Value = ""

# E8: Asked By, F8: Answered By
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = -1
